$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 first; this shifts old rows 13-21 down to 14-22,
# fixing the row-1-behind misalignment between labels (col A) and data (col B/C).
$ws.Rows("13:13").Insert()

# Row 10 (Objetivos:) - fill B/C with the correct Portuguese objectives text
# (previously these cells incorrectly held the "Messias Borges Silva" string)
$ws.Cells.Item(10, 2).Value = 'Ter uma maior conscientização entre os alunos sobre questões ligadas à área de sustentabilidadeCompreender o papel da engenharia e da tecnologia no desenvolvimento sustentável;Conhecer os métodos, ferramentas e incentivos para o desenvolvimento sustentável do sistema de produtos-serviçosEstabelecer uma compreensão clara do papel e do impacto de vários aspectos das decisões de engenharia sobre problemas ambientais, sociais e econômicos.'
$ws.Cells.Item(10, 3).Value = 'Ter uma maior conscientização entre os alunos sobre questões ligadas à área de sustentabilidadeCompreender o papel da engenharia e da tecnologia no desenvolvimento sustentável;Conhecer os métodos, ferramentas e incentivos para o desenvolvimento sustentável do sistema de produtos-serviçosEstabelecer uma compreensão clara do papel e do impacto de vários aspectos das decisões de engenharia sobre problemas ambientais, sociais e econômicos.'

# New row 13 holds the "Docentes responsaveis:" data (label stays on row 12)
$ws.Cells.Item(13, 2).Value = '5840535 - Messias Borges Silva'
$ws.Cells.Item(13, 3).Value = '5840535 - Messias Borges Silva'

# Row 14 (Programa resumido:) - new Portuguese short-syllabus text
$ws.Cells.Item(14, 2).Value = 'Sustentabilidade. Protocolos ambientais. Questões ambientais. Recursos naturais e sua poluição, créditos de carbono, conceito de resíduos zero, ISO 14000, análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial.'
$ws.Cells.Item(14, 3).Value = 'Sustentabilidade. Protocolos ambientais. Questões ambientais. Recursos naturais e sua poluição, créditos de carbono, conceito de resíduos zero, ISO 14000, análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial.'

# Row 16 (Programa:) - new Portuguese full syllabus text
$ws.Cells.Item(16, 2).Value = 'Sustentabilidade – necessidade, conceito, desafios, Protocolos ambientais,Questões ambientais globais, regionais e locais, Recursos naturais e sua poluição, Créditos de carbono, Conceito de resíduos zero, ISO 14000,Análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, , Materiais verdes, Energia, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial'
$ws.Cells.Item(16, 3).Value = 'Sustentabilidade – necessidade, conceito, desafios, Protocolos ambientais,Questões ambientais globais, regionais e locais, Recursos naturais e sua poluição, Créditos de carbono, Conceito de resíduos zero, ISO 14000,Análise do ciclo de vida, estudos de avaliação de impacto ambiental, habitat sustentável, , Materiais verdes, Energia, Fontes convencionais e renováveis, Tecnologia e desenvolvimento sustentável, Urbanização sustentável, Ecologia Industrial'

# Row 19 (Metodo:) - teaching method text
$ws.Cells.Item(19, 2).Value = 'Aulas Expositivas; trabalhos e seminários.'
$ws.Cells.Item(19, 3).Value = 'Aulas Expositivas; trabalhos e seminários.'

# Row 20 (Criterio:) - evaluation criteria text
$ws.Cells.Item(20, 2).Value = 'Avaliação dos trabalhos e apresentações ao longo do semestre'
$ws.Cells.Item(20, 3).Value = 'Avaliação dos trabalhos e apresentações ao longo do semestre'

# Row 21 (Norma de recuperacao:) - recovery norm formula text
$ws.Cells.Item(21, 2).Value = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação'
$ws.Cells.Item(21, 3).Value = 'NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota da recuperação'

# Row 22 (Bibliografia:) - bibliography text (label already shifted into place)
$ws.Cells.Item(22, 2).Value = 'ALLEN, D.T., SHONNARD, D.R. , Sustainable Engineering :concepts, design and case studies, Prentice Hall, 2015BLOKDIJK, G. , ISO14000 - Simple Steps to Win, Insights and Opportunities for Maxing out Success, Complete Publishing, 2015LAVE, L.B., HENDRICKSON, C.T. , Environmental Life Cycle Assessment of Goods and Services, Ed John Hopkins, 2006'
$ws.Cells.Item(22, 3).Value = 'ALLEN, D.T., SHONNARD, D.R. , Sustainable Engineering :concepts, design and case studies, Prentice Hall, 2015BLOKDIJK, G. , ISO14000 - Simple Steps to Win, Insights and Opportunities for Maxing out Success, Complete Publishing, 2015LAVE, L.B., HENDRICKSON, C.T. , Environmental Life Cycle Assessment of Goods and Services, Ed John Hopkins, 2006'

Write-Output "edit complete"
